# Apply the "Office Theme" colour palette to the deck's (single) slide
# master / theme part, replacing the current "Integral" palette.
#
# Theme colours aren't exposed through a document-import API in this
# host (ApplyTheme/ApplyTemplate need real .thmx files, which aren't
# available), so - exactly as real PowerPoint VBA automation does when
# it wants to recolour the active theme in place - each of the twelve
# theme colour slots is written individually via
# Slide.ThemeColorScheme.Colors(i).RGB. Editing through any slide
# updates the shared theme part used by every slide on that master.

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Office Theme colour scheme (msoThemeColorDark1 .. msoThemeColorFollowedHyperlink)
$officeTheme = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hyperlink
    "954F72"  # 12 followed hyperlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $tcs.Colors($i).RGB = HexToRgb $officeTheme[$i - 1]
}
